$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix A6: it was stored as text "79174445", make it a true number
$ws.Cells.Item(6, 1).Value = 79174445

# Append new redemption row 7 (phone column must stay text, not a number)
$ws.Cells.Item(7, 1).NumberFormat = "@"
$ws.Cells.Item(7, 1).Value = "79174445"
$ws.Cells.Item(7, 1).Style = "Normal"
$ws.Cells.Item(7, 2).Value = 20
$ws.Cells.Item(7, 3).Value = "2025-08-18T09:08:11"
